$wb = $excel.ActiveWorkbook

# Both the "展览" and "全部类型" sheets contain identical tables of
# convention listings; column F ("想去人数") values were refreshed.
$sheetNames = @("展览", "全部类型")

# Row -> new value for column F
$updates = @{
    10 = 158
    12 = 4920
    14 = 7195
    18 = 591
    35 = 1118
    37 = 1319
    45 = 2419
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
